{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per the supplied diff): append a new sentence\n// \" COVID-19 also problem for mental health.\" right after the existing\n// sentence that ends \"...have never been experienced before.\" inside the\n// \"Especially important now with pandemic...\" bullet under \"Thoughts\".\n//\n// (The remaining hunks in the diff only add <w:proofErr/> spell/grammar-\n// check markers and split existing runs at the same text boundaries -\n// those are artifacts Word's background proofing engine stamps into the\n// file and are not reachable/settable through the Word JS API, and they\n// do not change any visible text, so there is nothing else to apply.)\n\nconst body = context.document.body;\n\nconst searchText = \"have never been experienced before.\";\n\n// Locate the paragraph that contains the target sentence so we can append\n// the new sentence at the very end of that paragraph (this mirrors the\n// diff, which adds a brand-new <w:r> run at the end of the paragraph\n// rather than merging into the existing run).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(searchText) !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(`Could not find paragraph containing: \"${searchText}\"`);\n}\n\ntargetParagraph.insertText(\n  \" COVID-19 also problem for mental health.\",\n  Word.InsertLocation.end\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Target change (per the supplied diff): append a new sentence\n# \" COVID-19 also problem for mental health.\" right after the existing\n# sentence that ends \"...have never been experienced before.\" inside the\n# \"Especially important now with pandemic...\" bullet under \"Thoughts\".\n#\n# (The remaining hunks in the diff only add <w:proofErr/> spell/grammar-\n# check markers and split existing runs at the same text boundaries -\n# those are artifacts Word's background proofing engine stamps into the\n# file and are not reachable/settable through the Word object model, and\n# they do not change any visible text, so there is nothing else to apply.)\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"have never been experienced before.\")\n\nif (-not $found) {\n    throw \"Could not find target sentence: 'have never been experienced before.'\"\n}\n\n# $rng now spans exactly the matched text; collapse to its end point and\n# insert the new sentence there so it becomes a new run right after it.\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertAfter(\" COVID-19 also problem for mental health.\")\n"}
